$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "301.22"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.17%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "32.71"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.91%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.74%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07734"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.64%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.948"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-16.43%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.836"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.36%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.800"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.90%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9202"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.07%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1767"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.78%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07790"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.82%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08649"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-5.91%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03173"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "5.74%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1002"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.08%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001511"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.42%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005859"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.61%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.461"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.33%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.27%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.37%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.287"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.25%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "16.58%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.21%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001225"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.28%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004411"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.29%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.24%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01705"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-2.50%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04687"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.94%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007682"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "9.04%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.67%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002323"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "6.18%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01138"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "16.65%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006251"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.41%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.24%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.8204"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-28.83%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.24%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.24%"
